# Generate Report for Handoff
# A new handoff was recorded for 6622e7f0-b6c8-48c4-9847-2cb255f6e063, so its
# "Latest Handoff Date(time)" is refreshed on the Overview sheet (row 5) and
# on each locale sheet's row 5 ("Latest Handoff Datetime" column, E).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D5").Value = "2016-03-23 02:38:58"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E5").Value = "2016-03-23 02:38:55"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E5").Value = "2016-03-23 02:38:58"
